$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.822.34"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "1.629.58"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'215.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").Value = "'0.5069"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "'0.06470"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.81%  "

$ws.Range("D9").Value = "'0.2578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "'19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.26%  "

$ws.Range("D11").Value = "'0.07801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").Value = "'4.257"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.52%  "

$ws.Range("D13").Value = "1.627.28"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "1.853.37"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "'0.5570"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "

$ws.Range("D16").Value = "'63.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "

$ws.Range("D17").Value = "0.0₅7544"
$ws.Range("E17").Value = "  -2.61%  "

$ws.Range("D18").Value = "25.825.67"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "'194.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("D21").Value = "'4.297"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.23%  "

$ws.Range("D22").Value = "'9.812"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("D23").Value = "'6.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.57%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'1.819"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.01%  "

$ws.Range("D26").Value = "'140.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("D27").Value = "'0.1261"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "

$ws.Range("D28").Value = "'6.724"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("D29").Value = "'15.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").Value = "'0.04862"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").Value = "'3.278"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").Value = "'3.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").Value = "'1.553"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").Value = "'2.378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").Value = "'0.8937"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("D37").Value = "'2.571"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "1.133.33"
$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").Value = "'0.5464"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.38%  "

$ws.Range("D40").Value = "'0.01559"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "'5.564"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("D43").Value = "'0.7946"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "'97.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("D45").Value = "1.781.24"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -6.95%  "

$ws.Range("D47").Value = "'0.4438"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("D48").Value = "'55.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "

$ws.Range("D49").Value = "'0.05058"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").Value = "'7.616"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").Value = "'1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
